$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 55, pushing the existing weekly records (rows 55-66)
# down to rows 56-67, and fill the new row with this week's data.
$ws.Rows.Item(55).Insert()

$ws.Range("A55").Value = 11
$ws.Range("B55").Value = "Vega Monumental Concepción"
$ws.Range("C55").Value = "Bíobío"
$ws.Range("D55").Value = 44588
$ws.Range("E55").Value = 8
$ws.Range("F55").Value = 100112001
$ws.Range("G55").Value = "Berenjena"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 180
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 11000
$ws.Range("M55").Value = 10556
$ws.Range("N55").Value = "$/caja 60 unidades"
$ws.Range("O55").Value = "Región Metropolitana"
$ws.Range("P55").Value = 176
$ws.Range("Q55").Value = 60
$ws.Range("R55").Value = "Hortaliza"
